# Aggiornamento fino a 1/09/2021
# Appends new daily rows (358-366) to the worksheet, mirroring the style
# of the last existing data row (357: bold/bordered/centered date cell
# in column A, plain numeric cells in columns B-D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial (days since 1899-12-30), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newData = @(
    @(44432, 1, 2, 62.51953735542357),
    @(44433, 0, 2, 62.51953735542357),
    @(44434, 0, 2, 62.51953735542357),
    @(44435, 0, 1, 31.25976867771178),
    @(44436, 0, 1, 31.25976867771178),
    @(44437, 5, 6, 187.5586120662707),
    @(44438, 0, 6, 187.5586120662707),
    @(44439, 0, 5, 156.2988433885589),
    @(44440, 0, 5, 156.2988433885589)
)

$lastRow = 357
$destRow = $lastRow

foreach ($entry in $newData) {
    $destRow = $destRow + 1

    # Copy the previous row's formatting (keeps column A's date style,
    # plain style for B/C/D) down into the new row before writing values.
    $srcRange = $ws.Range("A" + ($destRow - 1) + ":D" + ($destRow - 1))
    $dstRange = $ws.Range("A" + $destRow + ":D" + $destRow)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($destRow, 1).Value = $entry[0]
    $ws.Cells.Item($destRow, 2).Value = $entry[1]
    $ws.Cells.Item($destRow, 3).Value = $entry[2]
    $ws.Cells.Item($destRow, 4).Value = $entry[3]
}

Write-Output ("Updated dimension through row " + $destRow)
